$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the five distinct server IP addresses (192.168.1.113-117) with a
# single shared value of 127.0.0.1 for rows 2-6 in column F ("IP").
$ws.Range("F2:F6").Value = "127.0.0.1"

# Update the active cell/selection left behind in the sheet view.
$ws.Range("F14").Select()
